$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17's B value was re-typed (integer, not the old float) - normalize it.
$ws.Range("B17").Value = 650

# New ranking rows appended below the existing table (rows 18-39).
$rows = @(
    @{ Row = 18; Name = 'Ana'; Score = 100 },
    @{ Row = 19; Name = 'Alexander'; Score = 400 },
    @{ Row = 20; Name = 'Jose'; Score = 700 },
    @{ Row = 21; Name = 'ola'; Score = 150 },
    @{ Row = 22; Name = 'lacrespa'; Score = 750 },
    @{ Row = 23; Name = 'Daddy Yankee'; Score = 750 },
    @{ Row = 24; Name = 'TUlip '; Score = 750 },
    @{ Row = 25; Name = 'Bryant Myers'; Score = 500 },
    @{ Row = 26; Name = 'Margarita'; Score = 500 },
    @{ Row = 27; Name = 'miku'; Score = 900 },
    @{ Row = 28; Name = 'Eduardo'; Score = 1300 },
    @{ Row = 29; Name = 'Jeff'; Score = 150 },
    @{ Row = 30; Name = 'Eduardo'; Score = 300 },
    @{ Row = 31; Name = 'ale'; Score = 900 },
    @{ Row = 32; Name = 'Elvi'; Score = 800 },
    @{ Row = 33; Name = 'camilo'; Score = 850 },
    @{ Row = 34; Name = 'Clubpinguin69'; Score = 850 },
    @{ Row = 35; Name = 'Jose :3'; Score = 1250 },
    @{ Row = 36; Name = 'Alan'; Score = 900 },
    @{ Row = 37; Name = 'crisd'; Score = 900 },
    @{ Row = 38; Name = 'Mery'; Score = 1300 },
    @{ Row = 39; Name = 'Alexander'; Score = 700 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Name
    $ws.Cells.Item($r.Row, 2).Value = $r.Score
}

# Selection / view state moved as the user scrolled/selected while editing.
$ws.Range("D29").Select()
